$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.405063152313232
$ws.Range("B1").Value = 1.723286509513855
$ws.Range("C1").Value = 3.213600635528564
$ws.Range("D1").Value = 1.477963805198669
$ws.Range("E1").Value = 0.8257070779800415
